$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 2.55  # G2: 2.6 -> 2.55
$ws.Cells.Item(2, 9).Value = 3.1  # I2: 3 -> 3.1
$ws.Cells.Item(2, 13).Value = 1.11  # M2: 1.13 -> 1.11
$ws.Cells.Item(2, 14).Value = 6.5  # N2: 6 -> 6.5
$ws.Cells.Item(2, 33).Value = 7  # AG2: 6.5 -> 7
$ws.Cells.Item(2, 37).Value = 34  # AK2: 29 -> 34
$ws.Cells.Item(2, 54).Value = 401  # BB2: 351 -> 401
$ws.Cells.Item(3, 13).Value = 1.1  # M3: 1.11 -> 1.1
$ws.Cells.Item(3, 14).Value = 7  # N3: 6.5 -> 7
$ws.Cells.Item(6, 7).Value = 3.3  # G6: 3.2 -> 3.3
$ws.Cells.Item(6, 9).Value = 2.15  # I6: 2.2 -> 2.15
$ws.Cells.Item(6, 41).Value = 19  # AO6: 17 -> 19
$ws.Cells.Item(6, 45).Value = 201  # AS6: 151 -> 201
$ws.Cells.Item(7, 7).Value = 1.8  # G7: 1.72 -> 1.8
$ws.Cells.Item(7, 8).Value = 3.3  # H7: 3.35 -> 3.3
$ws.Cells.Item(7, 9).Value = 4.35  # I7: 4.65 -> 4.35
$ws.Cells.Item(7, 10).Value = 2.35  # J7: 2.25 -> 2.35
$ws.Cells.Item(7, 11).Value = 2.1  # K7: 2.12 -> 2.1
$ws.Cells.Item(7, 12).Value = 4.6  # L7: 4.85 -> 4.6
$ws.Cells.Item(7, 14).Value = 7.8  # N7: 7.9 -> 7.8
$ws.Cells.Item(7, 15).Value = 1.33  # O7: 1.32 -> 1.33
$ws.Cells.Item(7, 16).Value = 2.8  # P7: 2.85 -> 2.8
$ws.Cells.Item(7, 17).Value = 1.98  # Q7: 1.93 -> 1.98
$ws.Cells.Item(7, 18).Value = 1.65  # R7: 1.7 -> 1.65
$ws.Cells.Item(7, 19).Value = 1.43  # S7: 1.42 -> 1.43
$ws.Cells.Item(7, 20).Value = 2.7  # T7: 2.74 -> 2.7
$ws.Cells.Item(7, 23).Value = 6.2  # W7: 6.3 -> 6.2
$ws.Cells.Item(7, 24).Value = 7.9  # X7: 7.8 -> 7.9
$ws.Cells.Item(7, 25).Value = 8.25  # Y7: 8 -> 8.25
$ws.Cells.Item(7, 26).Value = 14.5  # Z7: 13.5 -> 14.5
$ws.Cells.Item(7, 27).Value = 15.5  # AA7: 14.5 -> 15.5
$ws.Cells.Item(7, 28).Value = 30  # AB7: 28 -> 30
$ws.Cells.Item(7, 30).Value = 6.4  # AD7: 6.6 -> 6.4
$ws.Cells.Item(7, 31).Value = 15.5  # AE7: 16 -> 15.5
$ws.Cells.Item(7, 34).Value = 25  # AH7: 26 -> 25
$ws.Cells.Item(7, 35).Value = 14  # AI7: 15 -> 14
$ws.Cells.Item(7, 36).Value = 75  # AJ7: 80 -> 75
$ws.Cells.Item(7, 37).Value = 45  # AK7: 50 -> 45
$ws.Cells.Item(7, 38).Value = 50  # AL7: 55 -> 50
$ws.Cells.Item(7, 40).Value = 3.6  # AN7: 3.55 -> 3.6
$ws.Cells.Item(7, 41).Value = 8.75  # AO7: 8.25 -> 8.75
$ws.Cells.Item(7, 42).Value = 17.5  # AP7: 17 -> 17.5
$ws.Cells.Item(7, 43).Value = 30  # AQ7: 27 -> 30
$ws.Cells.Item(7, 44).Value = 60  # AR7: 55 -> 60
$ws.Cells.Item(7, 45).Value = 250  # AS7: 200 -> 250
$ws.Cells.Item(7, 46).Value = 2.7  # AT7: 2.72 -> 2.7
$ws.Cells.Item(7, 47).Value = 7.1  # AU7: 7.2 -> 7.1
$ws.Cells.Item(7, 48).Value = 60  # AV7: 65 -> 60
$ws.Cells.Item(7, 49).Value = 6.1  # AW7: 6.3 -> 6.1
$ws.Cells.Item(7, 50).Value = 24  # AX7: 26 -> 24
$ws.Cells.Item(7, 51).Value = 28  # AY7: 30 -> 28
$ws.Cells.Item(7, 53).Value = 150  # BA7: 175 -> 150
$ws.Cells.Item(7, 54).Value = 350  # BB7: 400 -> 350
$ws.Cells.Item(8, 7).Value = 2.72  # G8: 2.62 -> 2.72
$ws.Cells.Item(8, 8).Value = 3.05  # H8: 3 -> 3.05
$ws.Cells.Item(8, 9).Value = 2.52  # I8: 2.67 -> 2.52
$ws.Cells.Item(8, 12).Value = 3.05  # L8: 3.1 -> 3.05
$ws.Cells.Item(8, 14).Value = 8  # N8: 7.9 -> 8
$ws.Cells.Item(8, 17).Value = 1.95  # Q8: 1.98 -> 1.95
$ws.Cells.Item(8, 23).Value = 8.75  # W8: 8 -> 8.75
$ws.Cells.Item(8, 24).Value = 14.5  # X8: 13 -> 14.5
$ws.Cells.Item(8, 26).Value = 32  # Z8: 30 -> 32
$ws.Cells.Item(8, 28).Value = 30  # AB8: 32 -> 30
$ws.Cells.Item(8, 29).Value = 9  # AC8: 8.75 -> 9
$ws.Cells.Item(8, 30).Value = 5.9  # AD8: 5.8 -> 5.9
$ws.Cells.Item(8, 31).Value = 13  # AE8: 12.5 -> 13
$ws.Cells.Item(8, 32).Value = 60  # AF8: 55 -> 60
$ws.Cells.Item(8, 33).Value = 8  # AG8: 9 -> 8
$ws.Cells.Item(8, 34).Value = 12.5  # AH8: 14.5 -> 12.5
$ws.Cells.Item(8, 36).Value = 29  # AJ8: 32 -> 29
$ws.Cells.Item(8, 37).Value = 21  # AK8: 22 -> 21
$ws.Cells.Item(8, 38).Value = 30  # AL8: 28 -> 30
$ws.Cells.Item(8, 39).Value = 450  # AM8: 400 -> 450
$ws.Cells.Item(8, 40).Value = 4.7  # AN8: 4.6 -> 4.7
$ws.Cells.Item(8, 42).Value = 19.5  # AP8: 20 -> 19.5
$ws.Cells.Item(8, 44).Value = 80  # AR8: 90 -> 80
$ws.Cells.Item(8, 45).Value = 200  # AS8: 250 -> 200
$ws.Cells.Item(8, 46).Value = 2.62  # AT8: 2.6 -> 2.62
$ws.Cells.Item(8, 47).Value = 6.4  # AU8: 6.3 -> 6.4
$ws.Cells.Item(8, 49).Value = 4.5  # AW8: 4.7 -> 4.5
$ws.Cells.Item(8, 50).Value = 13  # AX8: 13.5 -> 13
$ws.Cells.Item(8, 51).Value = 19  # AY8: 18 -> 19
$ws.Cells.Item(8, 53).Value = 80  # BA8: 75 -> 80
$ws.Cells.Item(9, 12).Value = 3.6  # L9: 3.5 -> 3.6
$ws.Cells.Item(9, 19).Value = 1.33  # S9: 1.3 -> 1.33
$ws.Cells.Item(9, 20).Value = 3.25  # T9: 3.4 -> 3.25
$ws.Cells.Item(9, 21).Value = 1.57  # U9: 1.53 -> 1.57
$ws.Cells.Item(9, 22).Value = 2.25  # V9: 2.38 -> 2.25
$ws.Cells.Item(9, 27).Value = 17  # AA9: 15 -> 17
$ws.Cells.Item(9, 28).Value = 23  # AB9: 21 -> 23
$ws.Cells.Item(9, 29).Value = 13  # AC9: 15 -> 13
$ws.Cells.Item(9, 33).Value = 12  # AG9: 13 -> 12
$ws.Cells.Item(9, 35).Value = 12  # AI9: 11 -> 12
$ws.Cells.Item(9, 37).Value = 23  # AK9: 21 -> 23
$ws.Cells.Item(9, 38).Value = 29  # AL9: 26 -> 29
$ws.Cells.Item(9, 40).Value = 4.33  # AN9: 4.5 -> 4.33
$ws.Cells.Item(9, 41).Value = 11  # AO9: 12 -> 11
$ws.Cells.Item(9, 45).Value = 126  # AS9: 101 -> 126
$ws.Cells.Item(9, 46).Value = 3.25  # AT9: 3.4 -> 3.25
$ws.Cells.Item(9, 47).Value = 7.5  # AU9: 7 -> 7.5
$ws.Cells.Item(9, 50).Value = 17  # AX9: 15 -> 17
$ws.Cells.Item(9, 51).Value = 23  # AY9: 21 -> 23
$ws.Cells.Item(9, 53).Value = 67  # BA9: 51 -> 67
$ws.Cells.Item(9, 56).Value = 126  # BD9: 151 -> 126
$ws.Cells.Item(10, 8).Value = 3.6  # H10: 3.7 -> 3.6
$ws.Cells.Item(10, 9).Value = 3.1  # I10: 3 -> 3.1
$ws.Cells.Item(10, 37).Value = 23  # AK10: 21 -> 23
$ws.Cells.Item(10, 41).Value = 11  # AO10: 12 -> 11
$ws.Cells.Item(11, 13).Value = 1.11  # M11: 1.1 -> 1.11
$ws.Cells.Item(11, 14).Value = 6.5  # N11: 7 -> 6.5
$ws.Cells.Item(11, 17).Value = 2.6  # Q11: 2.5 -> 2.6
$ws.Cells.Item(11, 18).Value = 1.48  # R11: 1.5 -> 1.48
$ws.Cells.Item(11, 24).Value = 13  # X11: 15 -> 13
$ws.Cells.Item(11, 31).Value = 19  # AE11: 17 -> 19
$ws.Cells.Item(12, 7).Value = 3.6  # G12: 3.7 -> 3.6
$ws.Cells.Item(12, 12).Value = 2.75  # L12: 2.63 -> 2.75
$ws.Cells.Item(12, 23).Value = 8.5  # W12: 9 -> 8.5
$ws.Cells.Item(12, 24).Value = 17  # X12: 19 -> 17
$ws.Cells.Item(12, 40).Value = 5.5  # AN12: 6 -> 5.5
$ws.Cells.Item(13, 17).Value = 2.1  # Q13: 2.15 -> 2.1
$ws.Cells.Item(13, 18).Value = 1.7  # R13: 1.67 -> 1.7
$ws.Cells.Item(14, 17).Value = 2.05  # Q14: 2.08 -> 2.05
$ws.Cells.Item(14, 18).Value = 1.75  # R14: 1.73 -> 1.75
$ws.Cells.Item(16, 7).Value = 1.91  # G16: 1.95 -> 1.91
$ws.Cells.Item(16, 9).Value = 3.9  # I16: 3.75 -> 3.9
$ws.Cells.Item(16, 10).Value = 2.5  # J16: 2.6 -> 2.5
$ws.Cells.Item(16, 21).Value = 1.62  # U16: 1.57 -> 1.62
$ws.Cells.Item(16, 22).Value = 2.2  # V16: 2.25 -> 2.2
$ws.Cells.Item(16, 23).Value = 9  # W16: 9.5 -> 9
$ws.Cells.Item(16, 24).Value = 10  # X16: 11 -> 10
$ws.Cells.Item(16, 25).Value = 8.5  # Y16: 9 -> 8.5
$ws.Cells.Item(16, 37).Value = 29  # AK16: 26 -> 29
$ws.Cells.Item(16, 38).Value = 34  # AL16: 29 -> 34
$ws.Cells.Item(16, 39).Value = 151  # AM16: 126 -> 151
$ws.Cells.Item(16, 40).Value = 4  # AN16: 4.33 -> 4
$ws.Cells.Item(16, 48).Value = 51  # AV16: 41 -> 51
$ws.Cells.Item(18, 13).Value = 1.1  # M18: 1.08 -> 1.1
$ws.Cells.Item(18, 14).Value = 7  # N18: 8 -> 7
$ws.Cells.Item(18, 19).Value = 1.5  # S18: 1.53 -> 1.5
$ws.Cells.Item(18, 20).Value = 2.5  # T18: 2.38 -> 2.5
$ws.Cells.Item(18, 23).Value = 6  # W18: 5.5 -> 6
$ws.Cells.Item(18, 24).Value = 8  # X18: 7.5 -> 8
$ws.Cells.Item(18, 45).Value = 201  # AS18: 251 -> 201
$ws.Cells.Item(18, 46).Value = 2.5  # AT18: 2.38 -> 2.5
$ws.Cells.Item(18, 47).Value = 9  # AU18: 9.5 -> 9
$ws.Cells.Item(19, 14).Value = 9  # N19: 8.5 -> 9
$ws.Cells.Item(20, 15).Value = 1.5  # O20: 1.44 -> 1.5
$ws.Cells.Item(20, 16).Value = 2.5  # P20: 2.63 -> 2.5
$ws.Cells.Item(22, 7).Value = 3.25  # G22: 3.1 -> 3.25
$ws.Cells.Item(22, 9).Value = 2.07  # I22: 2.15 -> 2.07
$ws.Cells.Item(22, 10).Value = 3.8  # J22: 3.7 -> 3.8
$ws.Cells.Item(22, 11).Value = 2.07  # K22: 2.05 -> 2.07
$ws.Cells.Item(22, 12).Value = 2.7  # L22: 2.77 -> 2.7
$ws.Cells.Item(22, 14).Value = 6.45  # N22: 6.3 -> 6.45
$ws.Cells.Item(22, 16).Value = 2.67  # P22: 2.65 -> 2.67
$ws.Cells.Item(22, 18).Value = 1.62  # R22: 1.6 -> 1.62
$ws.Cells.Item(22, 20).Value = 2.47  # T22: 2.45 -> 2.47
$ws.Cells.Item(22, 23).Value = 8.75  # W22: 8.5 -> 8.75
$ws.Cells.Item(22, 24).Value = 16  # X22: 15 -> 16
$ws.Cells.Item(22, 25).Value = 11.75  # Y22: 11.5 -> 11.75
$ws.Cells.Item(22, 27).Value = 32  # AA22: 30 -> 32
$ws.Cells.Item(22, 33).Value = 6.6  # AG22: 6.7 -> 6.6
$ws.Cells.Item(22, 34).Value = 9.25  # AH22: 9.5 -> 9.25
$ws.Cells.Item(22, 35).Value = 9  # AI22: 9.25 -> 9
$ws.Cells.Item(22, 36).Value = 18.5  # AJ22: 19.5 -> 18.5
$ws.Cells.Item(22, 37).Value = 18.5  # AK22: 19 -> 18.5
$ws.Cells.Item(22, 40).Value = 5  # AN22: 4.85 -> 5
$ws.Cells.Item(22, 41).Value = 18  # AO22: 17.5 -> 18
$ws.Cells.Item(22, 46).Value = 2.45  # AT22: 2.4 -> 2.45
$ws.Cells.Item(22, 47).Value = 7.5  # AU22: 7.6 -> 7.5
$ws.Cells.Item(22, 48).Value = 75  # AV22: 80 -> 75
$ws.Cells.Item(22, 49).Value = 3.85  # AW22: 3.9 -> 3.85
$ws.Cells.Item(22, 50).Value = 10.75  # AX22: 11.25 -> 10.75
$ws.Cells.Item(22, 51).Value = 21  # AY22: 22 -> 21
$ws.Cells.Item(22, 52).Value = 40  # AZ22: 45 -> 40
$ws.Cells.Item(22, 53).Value = 80  # BA22: 90 -> 80
